$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is purely numeric-looking text (e.g. "225.67") need an
# explicit Text number format first, otherwise Excel auto-converts the assigned
# string into a real number (and e.g. "6.80" / "1.00" would lose trailing zeros).
$ws.Range('D2').Value = '95.708.81'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '3.623.68'
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('E4').Value = '  +27.71%  '
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '225.67'
$ws.Range('E6').Value = '  -5.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '639.56'
$ws.Range('E8').Value = '  -1.60%  '
$ws.Range('E9').Value = '  +3.10%  '
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '3.620.03'
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.83'
$ws.Range('E12').Value = '  +10.36%  '
$ws.Range('E13').Value = '  -0.30%  '
$ws.Range('E14').Value = '  -9.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.58'
$ws.Range('E15').Value = '  -3.40%  '
$ws.Range('D16').Value = '4.297.00'
$ws.Range('E16').Value = '  -2.23%  '
$ws.Range('D17').Value = '95.451.15'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '21.25'
$ws.Range('E18').Value = '  +13.73%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.82'
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.24'
$ws.Range('E20').Value = '  +9.63%  '
$ws.Range('D21').Value = '3.613.43'
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.519'
$ws.Range('E22').Value = '  +3.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '511.09'
$ws.Range('E23').Value = '  -1.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.25'
$ws.Range('E24').Value = '  -5.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.243'
$ws.Range('E25').Value = '  +26.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '123.25'
$ws.Range('E26').Value = '  +21.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000201'
$ws.Range('E27').Value = '  -8.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.80'
$ws.Range('E28').Value = '  -1.49%  '
$ws.Range('D29').Value = '3.814.54'
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.81'
$ws.Range('E30').Value = '  -5.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.17'
$ws.Range('E31').Value = '  +3.07%  '
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').Value = '  -4.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.620'
$ws.Range('E35').Value = '  +4.59%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.98'
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.78'
$ws.Range('E38').Value = '  -5.31%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '44.82'
$ws.Range('E39').Value = '  +11.64%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '596.19'
$ws.Range('E41').Value = '  -8.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.40'
$ws.Range('E42').Value = '  -5.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.05'
$ws.Range('E43').Value = '  +3.96%  '
$ws.Range('E44').Value = '  -2.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.483'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('E46').Value = '  +7.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.96'
$ws.Range('E47').Value = '  -3.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.948'
$ws.Range('E48').Value = '  -2.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.30'
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '225.31'
$ws.Range('E50').Value = '  +10.22%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.52'
$ws.Range('E51').Value = '  -0.38%  '
